$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: fill D8 across to X8 (same formula pattern, relative row refs) ---
$ws.Range("D8").Copy()
$ws.Range("E8:X8").PasteSpecial(-4122)   # xlPasteFormats: carry D8's style (s="7") to E8:X8
$ws.Range("D8:X8").FormulaR1C1 = "=(R[-2]C-R[-1]C)*R8C2+R9C2"

# --- Row 9: change D9's formula to reference $B$11 (pi_mult) instead of $A$13 ---
$ws.Range("D9").Formula = "=`$B`$11*EXP(D8)/(1+`$B`$11*EXP(D8))"
$ws.Range("D9").Copy()
$ws.Range("E9:X9").PasteSpecial(-4122)   # carry D9's style (s="7") to E9:X9
$ws.Range("D9:X9").FormulaR1C1 = "=R11C2*EXP(R[-1]C)/(1+R11C2*EXP(R[-1]C))"

# --- Row 10: fill D10 across to X10 (formula only, no formatting carried over) ---
$ws.Range("D10:X10").FormulaR1C1 = '=IF(OR(ISBLANK(R[-4]C), ISBLANK(R[-3]C)),"",CONCAT(ROUND(R[-1]C*100,0), "%"))'

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("O16").Select()
